# Delete the first data row (row 2: n-0 / school-0 / hali-0 / age 20).
# This shifts the remaining data rows (n-1..n-4) up by one, which also
# renumbers the "age" column values (21,22,23,24) onto rows 2-5 and
# shrinks the used range from A1:E6 down to A1:E5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Delete()
